# FK_Presentation Notes.docx - "Uploaded Final Second Presentation" edit
#
# Applies 5 changes:
#  1. Paragraph 4  - splits the "Ambinder, Psychologist & Head of Playtesting
#                     at Valve" run, adds a spellStart/spellEnd proofErr
#                     bracket around "Ambinder" and retitles Mike Ambinder.
#  2. Paragraph 5  - wraps "whether or not" in a gramStart/gramEnd proofErr
#                     bracket.
#  3. Paragraph 18 - wraps "Fiero" in a spellStart/spellEnd proofErr bracket.
#  4. Paragraph 22 - wraps the stray "to" in a spellStart/spellEnd proofErr
#                     bracket.
#  5. Paragraph 23 - removes the two text runs (the "Will flesh..." note and
#                     the trailing space run) while keeping the _GoBack
#                     bookmark.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. "...Mike Ambinder, Psychologist & Head of Playtesting at Valve" ---
$p = $d.Paragraphs(4)
$xml = @"
<w:p $wNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>“Playtesting is the most important part of the game development process</w:t></w:r><w:r><w:t xml:space="preserve">” – Mike </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ambinder</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>Senior Experimental Psychologist</w:t></w:r><w:r><w:t xml:space="preserve"> of Playtesting at Valve</w:t></w:r></w:p>
"@
$p.Range.InsertXML($xml)

# --- 2. "Playtesting [..] gains insight into whether or not the game is" ---
$p = $d.Paragraphs(5)
$xml = @"
<w:p $wNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>“</w:t></w:r><w:r><w:t xml:space="preserve">Playtesting [..] gains insight into </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>whether or not</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the game is </w:t></w:r><w:r><w:t>achieving</w:t></w:r><w:r><w:t xml:space="preserve"> your player experience goals”, Tracy Fullerton in “Game Design Workshop”</w:t></w:r></w:p>
"@
$p.Range.InsertXML($xml)

# --- 3. "Challenge; Fiero from achieving a difficult goal" ---
$p = $d.Paragraphs(18)
$xml = @"
<w:p $wNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Challenge; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Fiero</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from achieving a difficult goal</w:t></w:r></w:p>
"@
$p.Range.InsertXML($xml)

# --- 4. "...doesn't seem to make to much sense" ---
$p = $d.Paragraphs(22)
$xml = @"
<w:p $wNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">May </w:t></w:r><w:r><w:t xml:space="preserve">be worth relocating the current position on this slide to earlier on in the presentation; current position, on reflection, doesn’t seem to make </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>to</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> much sense</w:t></w:r></w:p>
"@
$p.Range.InsertXML($xml)

# --- 5. Remove the "Will flesh this area out..." runs, keep the bookmark ---
$p = $d.Paragraphs(23)
$xml = @"
<w:p $wNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@
$p.Range.InsertXML($xml)

Write-Output "Applied 5 edits"
